# Refresh the crypto price/volume table (columns D = Price, E = Volume(1h))
# for rows 2-51, as produced by the scheduled GitHub Actions scraper run.
#
# Several "Price" values look numeric (e.g. "1.0000", "0.9999") but must be
# stored as literal text (they are not real numbers - Excel would otherwise
# collapse "1.0000" to 1, or "28.10" to 28.1, losing the formatted digits).
# Setting NumberFormat to "@" (Text) before assigning those values forces
# Excel to keep them as-typed; values that already aren't number-like
# (e.g. "30.395.18", multiple dots) are left alone since Excel already
# stores them as text without any extra formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.395.18'
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").Value = '1.938.93'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7706'
$ws.Range("E5").Value = '  +8.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '248.01'
$ws.Range("E6").Value = '  -1.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.10'
$ws.Range("E8").Value = '  +1.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3214'
$ws.Range("E9").Value = '  -2.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07118'
$ws.Range("E10").Value = '  -2.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7867'
$ws.Range("E11").Value = '  -2.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08020'
$ws.Range("E12").Value = '  -0.61%  '

$ws.Range("D13").Value = '1.940.58'
$ws.Range("E13").Value = '  +0.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.387'
$ws.Range("E14").Value = '  -1.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '95.15'
$ws.Range("E15").Value = '  +0.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.60'
$ws.Range("E16").Value = '  -3.53%  '

$ws.Range("D17").Value = '30.397.46'
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '256.43'
$ws.Range("E18").Value = '  +1.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008023'
$ws.Range("E19").Value = '  -2.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.854'
$ws.Range("E20").Value = '  +1.00%  '

$ws.Range("D21").Value = '2.195.27'
$ws.Range("E21").Value = '  +0.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("E24").Value = '  -3.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.623'
$ws.Range("E25").Value = '  -1.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.29'
$ws.Range("E26").Value = '  -0.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.16'
$ws.Range("E27").Value = '  -1.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1345'
$ws.Range("E28").Value = '  +4.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.304'
$ws.Range("E29").Value = '  -1.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.366'
$ws.Range("E30").Value = '  +1.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.525'
$ws.Range("E31").Value = '  -1.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.445'
$ws.Range("E32").Value = '  +0.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.157'
$ws.Range("E33").Value = '  -0.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05205'
$ws.Range("E34").Value = '  +0.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.284'
$ws.Range("E35").Value = '  +1.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7521'
$ws.Range("E36").Value = '  +0.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.775'
$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01977'
$ws.Range("E38").Value = '  +0.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.809'
$ws.Range("E39").Value = '  -0.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.25'
$ws.Range("E40").Value = '  +0.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.470'
$ws.Range("E41").Value = '  +0.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4537'
$ws.Range("E42").Value = '  +0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.984'
$ws.Range("E43").Value = '  -1.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8366'
$ws.Range("E45").Value = '  -0.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.45'
$ws.Range("E46").Value = '  -0.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.806'
$ws.Range("E47").Value = '  +0.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.531'
$ws.Range("E48").Value = '  +1.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '985.83'
$ws.Range("E49").Value = '  +11.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.47'
$ws.Range("E50").Value = '  +2.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4173'
$ws.Range("E51").Value = '  -0.06%  '

